# The NATMI export was regenerated with updated TPM-derived stats.
# Net effect: the old 6-row table (2 ligand-expressing clusters x
# 3 target clusters minus blanks = 6 combos) collapses to 3 rows,
# one per sending cluster, all now pointing at "MuSCs" as the
# (single) target cluster, with refreshed numeric columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: ECs -> Cxcl1 -> Cxcr1 -> MuSCs ---------------------------------
$ws.Range("D2").Value = "MuSCs"
$ws.Range("G2").Value = 5.385511999999999
$ws.Range("H2").Value = 16.156536
$ws.Range("I2").Value = 0.02736372477514656
$ws.Range("J2").Value = 0.02736372477514657
$ws.Range("M2").Value = 0.000484
$ws.Range("N2").Value = 0.001452
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.002606587808
$ws.Range("R2").Value = 0.023459290272
$ws.Range("S2").Value = 0.02736372477514656
$ws.Range("T2").Value = 0.02736372477514657

# --- Row 3: FAPs -> Cxcl1 -> Cxcr1 -> MuSCs --------------------------------
$ws.Range("A3").Value = "FAPs"
$ws.Range("G3").Value = 153.7577133333333
$ws.Range("H3").Value = 461.27314
$ws.Range("I3").Value = 0.7812411799860843
$ws.Range("J3").Value = 0.7812411799860843
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.07441873325333333
$ws.Range("R3").Value = 0.66976859928
$ws.Range("S3").Value = 0.7812411799860843
$ws.Range("T3").Value = 0.7812411799860843

# --- Row 4: MuSCs -> Cxcl1 -> Cxcr1 -> MuSCs -------------------------------
$ws.Range("A4").Value = "MuSCs"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("G4").Value = 37.66886966666667
$ws.Range("H4").Value = 113.006609
$ws.Range("I4").Value = 0.1913950952387691
$ws.Range("J4").Value = 0.1913950952387691
$ws.Range("M4").Value = 0.000484
$ws.Range("N4").Value = 0.001452
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.01823173291866667
$ws.Range("R4").Value = 0.164085596268
$ws.Range("S4").Value = 0.1913950952387691
$ws.Range("T4").Value = 0.1913950952387691

# The former rows 5-7 (duplicate ECs-target / MuSCs-target combos) are gone;
# their useful data has already been folded into rows 2-4 above, so just
# delete the now-redundant trailing rows. This also shrinks the sheet
# dimension from A1:T7 to A1:T4 and drops the sharedStrings usage count.
$ws.Rows("5:7").Delete()
